$wb = $excel.ActiveWorkbook

$sheetsData = @(
    @{
        Name = "DE_LFT_#1"
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        D = "0x01,0x4C"
        E = "0x14"
        F = 380
        G = "7.598631275147109e+23" -as [double]
        H = 332
        I = 14
    },
    @{
        Name = "DE_LFT_#2"
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        D = "0x01,0x4C"
        E = "0xe"
        F = 380
        G = "5.68432987514711e+23" -as [double]
        H = 332
        I = 14
    },
    @{
        Name = "DE_PLT_#1"
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x00,0x7A"
        E = "0x7"
        F = 130
        G = "5.68631262647114e+23" -as [double]
        H = 122
        I = 7
    },
    @{
        Name = "DE_PLT_#2"
        B = "0x00,0x82"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x00,0x7A"
        E = "0x3"
        F = 130
        G = "9.85046333984776e+23" -as [double]
        H = 122
        I = 3
    }
)

foreach ($row in $sheetsData) {
    $ws = $wb.Worksheets.Item($row.Name)
    $newRow = 69

    $ws.Cells.Item($newRow, 1).Value = 45855.43355324074
    $ws.Cells.Item($newRow, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($newRow, 2).Value = $row.B
    $ws.Cells.Item($newRow, 3).Value = $row.C
    $ws.Cells.Item($newRow, 4).Value = $row.D
    $ws.Cells.Item($newRow, 5).Value = $row.E
    $ws.Cells.Item($newRow, 6).Value = $row.F
    $ws.Cells.Item($newRow, 7).Value = $row.G
    $ws.Cells.Item($newRow, 8).Value = $row.H
    $ws.Cells.Item($newRow, 9).Value = $row.I
}
